$d = $word.ActiveDocument

# Locate the paragraph that contains the astromap link with the old year (2018).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*GaNight/2018*") {
        $target = $p
        break
    }
}

$pr = $target.Range
$start = $pr.Start
$end = $pr.End - 1   # exclude the trailing paragraph mark

# Remove all the existing (heavily split-up) runs that make up the credit line.
$oldRange = $d.Range($start, $end)
$oldRange.Delete()

# Insert the new, consolidated credit line text (plain run, no explicit rPr).
$newTextRange = $d.Range($start, $start)
$newTextRange.InsertAfter("Jeník Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/")

# Insert a leading empty run ahead of the text run, matching the target markup.
$emptyRunRange = $d.Range($start, $start)
$emptyRunRange.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
